$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from existing styled header cells (A1:K1) to the newly added header cells (L1:W1)
$ws.Range("A1:K1").Copy()
$ws.Range("L1:W1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Header row values
$ws.Range("A1").Value = "productIds"
$ws.Range("B1").Value = "MSE_no_transfer"
$ws.Range("C1").Value = "MSE_transfer_basic"
$ws.Range("D1").Value = "MSE_transfer_coral"
$ws.Range("E1").Value = "MSE_transfer_sa"
$ws.Range("F1").Value = "MSE_transfer_bw"
$ws.Range("G1").Value = "MSE_transfer_nnw"
$ws.Range("H1").Value = "MAE_no_transfer"
$ws.Range("I1").Value = "MAE_transfer_basic"
$ws.Range("J1").Value = "MAE_transfer_coral"
$ws.Range("K1").Value = "MAE_transfer_sa"
$ws.Range("L1").Value = "MAE_transfer_bw"
$ws.Range("M1").Value = "MAE_transfer_nnw"
$ws.Range("N1").Value = "MSE_diff_basic"
$ws.Range("O1").Value = "MSE_transfer_coral"
$ws.Range("P1").Value = "MSE_diff_sa"
$ws.Range("Q1").Value = "MSE_diff_bw"
$ws.Range("R1").Value = "MSE_diff_nnw"
$ws.Range("S1").Value = "MAE_diff_basic"
$ws.Range("T1").Value = "MAE_transfer_coral"
$ws.Range("U1").Value = "MAE_diff_sa"
$ws.Range("V1").Value = "MAE_diff_bw"
$ws.Range("W1").Value = "MAE_diff_nnw"

# Data rows
# Row 2
$ws.Range("A2").Value = "101-120"
$ws.Range("B2").Value = 1.224559961485893
$ws.Range("C2").Value = 1.065808634844566
$ws.Range("D2").Value = 1.120835828974163
$ws.Range("E2").Value = 1.033296906461127
$ws.Range("F2").Value = 0.8809340079803796
$ws.Range("G2").Value = 0.890879775520842
$ws.Range("H2").Value = 0.5149422876031311
$ws.Range("I2").Value = 0.607803577699815
$ws.Range("J2").Value = 0.5656839662244296
$ws.Range("K2").Value = 0.78726599205821
$ws.Range("L2").Value = 0.4370141014276453
$ws.Range("M2").Value = 0.3903746227792839
$ws.Range("N2").Value = -0.1587513266413263
$ws.Range("O2").Value = -0.10372413251173
$ws.Range("P2").Value = -0.191263055024766
$ws.Range("Q2").Value = -0.3436259535055131
$ws.Range("R2").Value = -0.3336801859650507
$ws.Range("S2").Value = 0.09286129009668387
$ws.Range("T2").Value = 0.05074167862129852
$ws.Range("U2").Value = 0.2723237044550789
$ws.Range("V2").Value = -0.0779281861754858
$ws.Range("W2").Value = -0.1245676648238472

# Row 3
$ws.Range("A3").Value = "121-140"
$ws.Range("B3").Value = 2.588536900271843
$ws.Range("C3").Value = 2.592797257199938
$ws.Range("D3").Value = 3.161065882811345
$ws.Range("E3").Value = 2.852513077306079
$ws.Range("F3").Value = 2.806361749175249
$ws.Range("G3").Value = 2.699051178012179
$ws.Range("H3").Value = 0.9914813140809273
$ws.Range("I3").Value = 0.9085127110246211
$ws.Range("J3").Value = 0.8625860071605033
$ws.Range("K3").Value = 1.032059815070777
$ws.Range("L3").Value = 0.8616075054061442
$ws.Range("M3").Value = 0.7461395487668154
$ws.Range("N3").Value = 0.004260356928094566
$ws.Range("O3").Value = 0.5725289825395015
$ws.Range("P3").Value = 0.2639761770342361
$ws.Range("Q3").Value = 0.2178248489034056
$ws.Range("R3").Value = 0.1105142777403354
$ws.Range("S3").Value = -0.08296860305630627
$ws.Range("T3").Value = -0.1288953069204241
$ws.Range("U3").Value = 0.04057850098984983
$ws.Range("V3").Value = -0.1298738086747832
$ws.Range("W3").Value = -0.2453417653141119

# Row 4
$ws.Range("A4").Value = "141-160"
$ws.Range("B4").Value = 5.485710658411413
$ws.Range("C4").Value = 5.202870679331181
$ws.Range("D4").Value = 3.502699856861465
$ws.Range("E4").Value = 3.176874955584172
$ws.Range("F4").Value = 3.230494342152407
$ws.Range("G4").Value = 3.11746815953909
$ws.Range("H4").Value = 1.282621473071498
$ws.Range("I4").Value = 1.286737497242299
$ws.Range("J4").Value = 1.009296810350896
$ws.Range("K4").Value = 1.108763838657743
$ws.Range("L4").Value = 0.9005700749771435
$ws.Range("M4").Value = 0.8744593543504965
$ws.Range("N4").Value = -0.2828399790802321
$ws.Range("O4").Value = -1.983010801549948
$ws.Range("P4").Value = -2.30883570282724
$ws.Range("Q4").Value = -2.255216316259006
$ws.Range("R4").Value = -2.368242498872323
$ws.Range("S4").Value = 0.004116024170801591
$ws.Range("T4").Value = -0.2733246627206014
$ws.Range("U4").Value = -0.1738576344137543
$ws.Range("V4").Value = -0.3820513980943543
$ws.Range("W4").Value = -0.4081621187210013

# Row 5
$ws.Range("A5").Value = "161-180"
$ws.Range("B5").Value = 14.22789284129596
$ws.Range("C5").Value = 13.37611582038006
$ws.Range("D5").Value = 12.53143968824912
$ws.Range("E5").Value = 12.63280682930526
$ws.Range("F5").Value = 12.37904794826922
$ws.Range("G5").Value = 12.56539049160015
$ws.Range("H5").Value = 1.136556738826326
$ws.Range("I5").Value = 1.163588824755907
$ws.Range("J5").Value = 1.015973113539503
$ws.Range("K5").Value = 1.018982735495265
$ws.Range("L5").Value = 0.9943578763063982
$ws.Range("M5").Value = 0.9724489795748279
$ws.Range("N5").Value = -0.8517770209158986
$ws.Range("O5").Value = -1.696453153046841
$ws.Range("P5").Value = -1.595086011990695
$ws.Range("Q5").Value = -1.848844893026735
$ws.Range("R5").Value = -1.662502349695812
$ws.Range("S5").Value = 0.02703208592958162
$ws.Range("T5").Value = -0.1205836252868224
$ws.Range("U5").Value = -0.1175740033310602
$ws.Range("V5").Value = -0.1421988625199274
$ws.Range("W5").Value = -0.1641077592514977

# Row 6
$ws.Range("A6").Value = "181-200"
$ws.Range("B6").Value = 8.201373654786686
$ws.Range("C6").Value = 7.638955692436562
$ws.Range("D6").Value = 5.859342160701627
$ws.Range("E6").Value = 5.807247552767949
$ws.Range("F6").Value = 5.88428453792401
$ws.Range("G6").Value = 6.225046535432612
$ws.Range("H6").Value = 1.062188837815168
$ws.Range("I6").Value = 1.063140948999174
$ws.Range("J6").Value = 0.8617176915191481
$ws.Range("K6").Value = 0.8229994240953901
$ws.Range("L6").Value = 0.808653786169095
$ws.Range("M6").Value = 0.7944319410373135
$ws.Range("N6").Value = -0.562417962350124
$ws.Range("O6").Value = -2.342031494085059
$ws.Range("P6").Value = -2.394126102018737
$ws.Range("Q6").Value = -2.317089116862675
$ws.Range("R6").Value = -1.976327119354074
$ws.Range("S6").Value = 0.0009521111840060659
$ws.Range("T6").Value = -0.2004711462960199
$ws.Range("U6").Value = -0.239189413719778
$ws.Range("V6").Value = -0.253535051646073
$ws.Range("W6").Value = -0.2677568967778545
